$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (price) changes - must remain as text, matching original inlineStr cells.
# Force text via NumberFormat "@" before assignment, then restore default "Normal" style
# so no stray style index is left attached to the cell.
$dValues = [ordered]@{
    'D2' = '37.196.86'
    'D3' = '2.064.76'
    'D5' = '250.88'
    'D6' = '0.678'
    'D7' = '60.97'
    'D9' = '61.23'
    'D11' = '0.0802'
    'D13' = '15.37'
    'D14' = '2.363.00'
    'D15' = '0.820'
    'D16' = '5.38'
    'D17' = '2.053.79'
    'D18' = '37.132.12'
    'D19' = '75.33'
    'D20' = '0.0₃0931'
    'D21' = '14.58'
    'D22' = '5.40'
    'D23' = '240.08'
    'D25' = '2.44'
    'D26' = '171.60'
    'D27' = '9.24'
    'D28' = '20.37'
    'D31' = '4.64'
    'D33' = '0.0636'
    'D34' = '4.42'
    'D36' = '0.999'
    'D37' = '2.29'
    'D39' = '0.113'
    'D41' = '18.41'
    'D43' = '1.15'
    'D44' = '97.78'
    'D45' = '4.38'
    'D47' = '4.60'
    'D48' = '2.53'
    'D49' = '1.307.25'
    'D50' = '2.93'
    'D51' = '6.92'
}
foreach ($addr in $dValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$addr]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Column E (volume/percentage) changes - plain text already (contains spaces/%).
$eValues = [ordered]@{
    'E2' = '  +0.27%  '
    'E3' = '  -1.00%  '
    'E4' = '  +0.21%  '
    'E5' = '  +0.03%  '
    'E6' = '  +3.47%  '
    'E7' = '  +22.46%  '
    'E8' = '  -0.01%  '
    'E9' = '  +1.21%  '
    'E10' = '  +2.89%  '
    'E11' = '  +8.00%  '
    'E13' = '  +1.91%  '
    'E14' = '  -1.25%  '
    'E15' = '  -0.86%  '
    'E16' = '  +5.53%  '
    'E17' = '  -2.69%  '
    'E18' = '  +0.65%  '
    'E19' = '  +4.59%  '
    'E20' = '  +13.23%  '
    'E21' = '  +10.31%  '
    'E22' = '  +4.04%  '
    'E23' = '  +0.11%  '
    'E24' = '  +0.14%  '
    'E25' = '  -0.52%  '
    'E26' = '  +1.83%  '
    'E27' = '  -0.52%  '
    'E28' = '  -2.24%  '
    'E29' = '  +0.87%  '
    'E30' = '  +2.72%  '
    'E31' = '  +3.71%  '
    'E32' = '  -4.23%  '
    'E33' = '  +5.31%  '
    'E34' = '  +8.38%  '
    'E35' = '  -2.98%  '
    'E37' = '  +0.52%  '
    'E38' = '  -3.24%  '
    'E39' = '  +29.75%  '
    'E40' = '  +3.14%  '
    'E41' = '  +4.57%  '
    'E42' = '  +1.59%  '
    'E43' = '  +0.60%  '
    'E44' = '  +0.44%  '
    'E45' = '  +32.41%  '
    'E46' = '  -0.05%  '
    'E47' = '  +17.49%  '
    'E48' = '  +12.45%  '
    'E49' = '  +0.08%  '
    'E50' = '  -0.58%  '
    'E51' = '  +0.71%  '
}
foreach ($addr in $eValues.Keys) {
    $ws.Range($addr).Value = $eValues[$addr]
}

# Column B (coin name) changes for the reshuffled rows 43-45.
$bValues = [ordered]@{
    'B43' = 'ARBITRUM'
    'B44' = 'Aave'
    'B45' = 'FTXToken'
}
foreach ($addr in $bValues.Keys) {
    $ws.Range($addr).Value = $bValues[$addr]
}

# Column C (link) changes for the reshuffled rows 43-45.
$cValues = [ordered]@{
    'C43' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'C44' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'C45' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
}
foreach ($addr in $cValues.Keys) {
    $ws.Range($addr).Value = $cValues[$addr]
}
